# Fixed popup menu style
# Adds a new localization row ("more_options") to the "string" table/sheet,
# and widens column E (lu) to fit the longer localized text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Append the new translation row (row 96) ---
$ws.Range("A96").Value = "more_options"
$ws.Range("B96").Value = "More options"
$ws.Range("C96").Value = "Plus d'options"
$ws.Range("D96").Value = "Mehr Optionen"
$ws.Range("E96").Value = "Méi Optiounen"

# --- Grow the "string" table so the new row is included ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E96"))

# --- Update the ExternalData_1 defined name so it spans the extra row ---
$wb.Names.Item(1).RefersTo = "=string!`$A`$1:`$B`$96"

# --- Widen column E (lu) to better fit the localized text, no more autofit bestFit ---
$ws.Columns.Item(5).ColumnWidth = 32.8

# --- Update selection / active cell to the newly added cell ---
$ws.Range("E96").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 76
$win.ScrollColumn = 1
